# Meeting Report.xlsx -- "Refactor of IMU test and synced share"
#
# 1. Fill in the Week-18 meeting row (date + hours grid) that was logged
#    after the fact -- all the downstream subtotal/summary formulas on the
#    sheet recalc automatically from this.
# 2. Collapse/hide the "Pre PDR Hours" detail columns (X:AC), mirroring the
#    already-collapsed "Pre PDR Costs" group (P:U).
# 3. Move the frozen-pane selection to J19.
# 4. Switch the workbook to manual calculation before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Week of 3/19/18 (row 18): fill in attendance + hours -----------
$ws.Range("A18").Value = 43202
$ws.Range("C18").Value = 0.5
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 1

# --- 2. Collapse the Pre-PDR-Hours detail columns (X:AC) ----------------
$detailCols = $ws.Range("X1:AC1").EntireColumn
$detailCols.OutlineLevel = 1
$detailCols.Hidden = $true

# --- 3. Update the active selection -------------------------------------
$ws.Range("J19").Select()

# --- 4. Save with manual calculation (values already fresh) -------------
$excel.Calculation = -4135
